# "updated forms with labs"
# The survey's "region" question is relabelled to refer to a facility lab,
# and the quick-search it drives switches from the 'regions' list to 'labs'.
# The settings sheet (previously just form_title/form_id/version metadata)
# is rebuilt into a choices-style "list name / name / label" lookup table,
# adding a 'lab' list (keyed lab_key/lab) alongside the existing
# facility/stype/condition lookups already present on the choices sheet.

$wb = $excel.ActiveWorkbook

# --- survey sheet -----------------------------------------------------
$survey = $wb.Worksheets.Item("survey")
$survey.Range("C6").Value = "Facility Lab"
$survey.Range("E6").Value = "quick search('labs')"

# --- settings sheet -----------------------------------------------------
$settings = $wb.Worksheets.Item("settings")

# grow the sheet from 2 rows to 5, keeping row 1's header formatting and
# row 2's data formatting for the newly inserted rows 3-5
$settings.Range("A3:A5").EntireRow.Insert()

$settings.Range("A1").Value = "list name"
$settings.Range("B1").Value = "name"
$settings.Range("C1").Value = "label"

$settings.Range("A2").Value = "region"
$settings.Range("B2").Value = "lab_key"
$settings.Range("C2").Value = "lab"

$settings.Range("A3").Value = "facility"
$settings.Range("B3").Value = "facility_key"
$settings.Range("C3").Value = "facility"

$settings.Range("A4").Value = "stype"
$settings.Range("B4").Value = "stype_key"
$settings.Range("C4").Value = "stype"

$settings.Range("A5").Value = "condition"
$settings.Range("B5").Value = "cond_key"
$settings.Range("C5").Value = "cond"

Write-Host "updated forms with labs: done"
